$wb = $excel.ActiveWorkbook

$wsContents = $wb.Worksheets.Item("Contents")
$wsPop      = $wb.Worksheets.Item("Population control")
$wsIncome   = $wb.Worksheets.Item("Income control")
$wsPreExc   = $wb.Worksheets.Item("Pre-tax top shares exc. CG")
$wsPostExc  = $wb.Worksheets.Item("Post-tax top shares exc. CG")
$wsPreInc   = $wb.Worksheets.Item("Pre-tax top shares inc. CG")

# ---------------------------------------------------------------------------
# Fix the mislabelled title on the "Pre-tax top shares inc. CG" sheet: it was
# pointing at the "excluding capital gains" string; give it its own text.
# ---------------------------------------------------------------------------
$wsPreInc.Range("A1").Value = "Pre-tax top shares including realised capital gains (%)"

# ---------------------------------------------------------------------------
# Updated underlying data (revised population / income control figures).
# ---------------------------------------------------------------------------

$sheet3Data = @{
  10 = @{ "B" = 665033.58865569381; "C" = 521047.59501383384; "D" = 680775.18465958792 }
  11 = @{ "B" = 729672.15384509147; "C" = 570420.59404806374; "D" = 750576.44599025627 }
  12 = @{ "B" = 796668.35196319094; "C" = 619348.60864261515; "D" = 824639.48969494947 }
  13 = @{ "B" = 849287.20391328097; "C" = 658666.03601612139; "D" = 888502.28042037482 }
  14 = @{ "B" = 912828.67688211065; "C" = 706534.98697444552; "D" = 979445.85072795476 }
  18 = @{ "B" = 951249.36885674007; "C" = 748113.13546705153; "D" = 971543.70734174002 }
  23 = @{ "B" = 1141089.8431397646; "C" = 913595.30168744014; "D" = 1188949.2458607645 }
  24 = @{ "B" = 1178067.578536235;  "C" = 951723.70502621564; "D" = 1232781.428982235 }
}

$sheet4Data = @{
  10 = @{ "B" = 41.105871029390073; "C" = 28.720356046286852; "D" = 13.143485686821215; "E" = 9.5251284871937223; "F" = 4.5673545494966392; "G" = 3.2922538468922409 }
  11 = @{ "B" = 40.610897302647615; "C" = 28.400297682316857; "D" = 13.219970723537857; "E" = 9.6260889550227819; "F" = 4.6753289099711459; "G" = 3.4002329242568372 }
  12 = @{ "B" = 41.411333736282671; "C" = 29.402217227387951; "D" = 14.143643413460353; "E" = 10.405043128688462; "F" = 5.1469211445218335; "G" = 3.7632974252686182 }
  13 = @{ "B" = 41.8397971468912;   "C" = 29.990269820678147; "D" = 14.79414515882843;  "E" = 10.998315100565337; "F" = 5.5632659826938138; "G" = 4.0884711912246452 }
  14 = @{ "B" = 41.986874338476888; "C" = 30.316727487030754; "D" = 15.223292316837419; "E" = 11.434140690035379; "F" = 5.9667096717103272; "G" = 4.4486590958981553 }
  18 = @{ "B" = 40.961532004171708; "C" = 28.901987802791588; "D" = 13.553895937748297; "E" = 9.9507685909555086; "F" = 5.0336203231234062; "G" = 3.7380995495599638 }
  23 = @{ "B" = 40.468691790261254; "C" = 28.98947997099085;  "D" = 14.11397906994924;  "E" = 10.617925465040345; "F" = 5.5945631981359929; "G" = 4.2101548648432487 }
}

$sheet5Data = @{
  10 = @{ "B" = 37.05619545982853;  "C" = 25.075187329246791; "D" = 10.907310167634437; "E" = 7.8232677592597168; "F" = 3.7249146974482827; "G" = 2.6916042600454086 }
  11 = @{ "B" = 36.602706636411412; "C" = 24.761907686231176; "D" = 10.949642129612661; "E" = 7.8830271918361916; "F" = 3.7846525533989848; "G" = 2.7475129565273977 }
  12 = @{ "B" = 37.353108421949216; "C" = 25.649922267925358; "D" = 11.763563652240187; "E" = 8.5675479226732758; "F" = 4.1963370229460013; "G" = 3.0657975937691817 }
  13 = @{ "B" = 37.818267053867366; "C" = 26.242494564937818; "D" = 12.406768287118989; "E" = 9.1455081037097692; "F" = 4.5873176674175831; "G" = 3.3622106486148109 }
  14 = @{ "B" = 37.871665119284167; "C" = 26.447058715680335; "D" = 12.731728791798194; "E" = 9.4798504936275894; "F" = 4.9034425321942976; "G" = 3.6428203186851529 }
  18 = @{ "B" = 36.026671375241378; "C" = 24.251591922312198; "D" = 10.371554565755352; "E" = 7.4062778246022321; "F" = 3.6053922429215284; "G" = 2.6580487576500098 }
  23 = @{ "B" = 35.271301173611981; "C" = 24.002738625851638; "D" = 10.714686448402064; "E" = 7.9105562767786886; "F" = 4.1136240130555146; "G" = 3.104980243844131 }
  24 = @{ "B" = 35.175148455512847; "C" = 24.010956597880977; "D" = 10.928029616428191; "E" = 8.1105226199284584; "F" = 4.2733644653872629; "G" = 3.2469681099663132 }
}

$sheet6Data = @{
  10 = @{ "B" = 14.338212462565112; "C" = 5.6630879771346372; "D" = 2.2232281496327553 }
  11 = @{ "B" = 14.761596052271628; "C" = 5.9921235044303174; "D" = 2.4010724748152996 }
  12 = @{ "B" = 15.471801883971761; "C" = 6.7821733759780054; "D" = 2.9433208082046911 }
  13 = @{ "B" = 16.906682945069651; "C" = 7.9731102832788032; "D" = 3.7228258117042183 }
  14 = @{ "B" = 18.917874203818943; "C" = 9.861528626721805;  "D" = 4.9688713725168627 }
  18 = @{ "B" = 14.379178100204667; "C" = 5.9554706147304364; "D" = 2.3941280072352225 }
  23 = @{ "B" = 16.426268882370039; "C" = 7.7076460848970321; "D" = 3.425713935375994 }
  24 = @{ "B" = 17.245555051516245; "C" = 8.2658610524435812; "D" = 3.6811067179577011 }
}

foreach ($row in $sheet3Data.Keys) {
  foreach ($col in $sheet3Data[$row].Keys) {
    $wsIncome.Range("$col$row").Value = $sheet3Data[$row][$col]
  }
}

foreach ($row in $sheet4Data.Keys) {
  foreach ($col in $sheet4Data[$row].Keys) {
    $wsPreExc.Range("$col$row").Value = $sheet4Data[$row][$col]
  }
}

foreach ($row in $sheet5Data.Keys) {
  foreach ($col in $sheet5Data[$row].Keys) {
    $wsPostExc.Range("$col$row").Value = $sheet5Data[$row][$col]
  }
}

foreach ($row in $sheet6Data.Keys) {
  foreach ($col in $sheet6Data[$row].Keys) {
    $wsPreInc.Range("$col$row").Value = $sheet6Data[$row][$col]
  }
}

# ---------------------------------------------------------------------------
# Restore each sheet's last-used selection, finishing on "Income control" so
# it ends up as the active (tabSelected) sheet, matching the saved view.
# ---------------------------------------------------------------------------
$wsContents.Activate()
$wsContents.Range("A4").Select() | Out-Null

$wsPop.Activate()
$wsPop.Range("E15").Select() | Out-Null

$wsPreExc.Activate()
$wsPreExc.Range("B3:G24").Select() | Out-Null

$wsPostExc.Activate()
$wsPostExc.Range("B3:G24").Select() | Out-Null

$wsPreInc.Activate()
$wsPreInc.Range("B2").Select() | Out-Null

$wsIncome.Activate()
$wsIncome.Range("F13").Select() | Out-Null
